$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1, reusing the same formatting (style) as the
# existing header cells (e.g. G1: bold, centered, bordered).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Fill the new "Save" column (H2:H9) with 0 for every data row.
for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 8).Value = 0
}
